$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ58519492",
    "summ59393415",
    "summ00177734",
    "summ01046587",
    "summ01923931",
    "summ02995260",
    "summ05115482",
    "summ07220522",
    "summ09218911"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
